# Nexial unit-test workbook: add a new "localdb" command type (and its
# sub-commands) to the '#system' sheet, shifting the existing command
# columns one slot to the right to make room for it at column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1. Insert a new column before N -- everything from N..AC slides to O..AD.
$ws.Columns("N:N").Insert(-4121)

# 2. Populate the new "localdb" command-family column (header + its verbs).
$localdb = @(
    "localdb",
    "cloneTable(var,source,target)",
    "dropTables(var,tables)",
    "exportCSV(sql,output)",
    "importRecords(var,sourceDb,sql,table)",
    "purge(var)",
    "runSQLs(var,sqls)"
)
for ($i = 0; $i -lt $localdb.Length; $i++) {
    $ws.Cells.Item(1 + $i, 14).Value = $localdb[$i]
}

# 3. Column A is the sorted index of all command-family names ("target").
#    "localdb" sorts in right after "json"/before "macro", so rows 14-29
#    shift down to 15-30 and row 14 becomes "localdb".
$colA = @("localdb","macro","mail","number","pdf","rdbms","redis","sms","sound","ssh","step","web","webalert","webcookie","ws","ws.async","xml")
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item(14 + $i, 1).Value = $colA[$i]
}

# 4. Fix up the defined names that pointed at the now-shifted columns
#    (and extend "target" to cover the new row 30).
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

# 5. New defined name for the "localdb" command-family range.
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

Write-Output "done"
